# Insert a new data row at row 151 (pushing existing rows 151-263 down to
# 152-264) and populate it with the new weekly record. All other rows shift
# down by one but keep their original values untouched by this script.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(151).Insert()

$ws.Range("A151").Value = 4
$ws.Range("B151").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C151").Value = 'Los Lagos'
$ws.Range("D151").Value = 44827
$ws.Range("E151").Value = 10
$ws.Range("F151").Value = 100112039
$ws.Range("G151").Value = 'Ciboulette'
$ws.Range("H151").Value = 'Sin especificar'
$ws.Range("I151").Value = 'Primera'
$ws.Range("J151").Value = 240
$ws.Range("K151").Value = 3000
$ws.Range("L151").Value = 3000
$ws.Range("M151").Value = 3000
$ws.Range("N151").Value = '$/docena de atados'
$ws.Range("O151").Value = 'Región Metropolitana'
$ws.Range("P151").Value = 1000
$ws.Range("Q151").Value = 3
$ws.Range("R151").Value = 'Hortaliza'
